# Update the crypto price/volume table with the latest scraped values.
# Rows 39/40, 42/43 and 44/45 also swapped rank order (Coin/Link/Price/Volume
# all move together), so those rows set B/C/D/E explicitly rather than just D/E.
# A leading apostrophe is used on numeric-looking Price values (single '.')
# so Excel stores them as text, matching the original text-typed cells
# (values with two dots, e.g. "28.347.17", are already non-numeric text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.347.17'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.811.84'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').Value = '''0.9999'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '''313.55'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('D8').Value = '''0.3993'
$ws.Range('E8').Value = '  +2.96%  '
$ws.Range('D9').Value = '''0.07879'
$ws.Range('E9').Value = '  -5.02%  '
$ws.Range('D10').Value = '''1.115'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').Value = '''40.96'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').Value = '''6.387'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '''0.9999'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').Value = '''20.42'
$ws.Range('E14').Value = '  -3.89%  '
$ws.Range('D15').Value = '''7.355'
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').Value = '1.804.04'
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '''92.80'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').Value = '''0.00001084'
$ws.Range('E18').Value = '  -3.62%  '
$ws.Range('D19').Value = '''0.06577'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').Value = '''0.9996'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '''17.36'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').Value = '''6.026'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').Value = '28.398.41'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').Value = '''11.20'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').Value = '''2.240'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = '''161.12'
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('D27').Value = '''20.55'
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').Value = '2.018.31'
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').Value = '''2.403'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '''128.91'
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').Value = '''0.1088'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').Value = '''1.070'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('D33').Value = '''3.663'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = '''5.592'
$ws.Range('E34').Value = '  -2.65%  '
$ws.Range('D35').Value = '''0.07263'
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('D36').Value = '''9.156'
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('D37').Value = '''0.02347'
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('D38').Value = '''0.2181'
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '''11.66'
$ws.Range('E39').Value = '  -2.84%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '''5.084'
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('D41').Value = '''0.6212'
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '''0.9996'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''1.160'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''13.26'
$ws.Range('E44').Value = '  -3.48%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.6009'
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('E46').Value = '  -6.05%  '
$ws.Range('D47').Value = '''3.744'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').Value = '''125.75'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').Value = '''1.225'
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').Value = '''0.06853'
$ws.Range('E51').Value = '  -1.85%  '
